$d = $word.ActiveDocument
$d.Content.Find.Execute("caretteristiche", $true, $false, $false, $false, $false,
                         $true, 1, $false, "caratteristiche", 2)
